$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 3015
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 3015
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 9045
$ws.Cells.Item(69, 13).Value = ""
$ws.Cells.Item(69, 14).Value = -10793

$ws.Cells.Item(72, 8).Value = 3015
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 3015
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 27135
$ws.Cells.Item(72, 13).Value = ""
$ws.Cells.Item(72, 14).Value = -35871

$ws.Cells.Item(112, 8).Value = 25001390
$ws.Cells.Item(112, 10).Value = 1502.5278
$ws.Cells.Item(112, 12).Value = 4507.5834
$ws.Cells.Item(112, 14).Value = -6723.5834

$ws.Cells.Item(116, 8).Value = 1684.6154
$ws.Cells.Item(116, 9).Value = 1587.5
$ws.Cells.Item(116, 10).Value = 1840
$ws.Cells.Item(116, 11).Value = 1587.5
$ws.Cells.Item(116, 12).Value = 1840
$ws.Cells.Item(116, 13).Value = 1854.5
$ws.Cells.Item(116, 14).Value = -8724

$ws.Cells.Item(132, 8).Value = 1043789.56
$ws.Cells.Item(132, 9).Value = 1213.907
$ws.Cells.Item(132, 10).Value = 12251478
$ws.Cells.Item(132, 11).Value = 3641.721
$ws.Cells.Item(132, 12).Value = 36754434
$ws.Cells.Item(132, 13).Value = -1111.721
$ws.Cells.Item(132, 14).Value = -36759494

$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).Value = ""

$ws.Cells.Item(140, 8).Value = 64550.24
$ws.Cells.Item(140, 10).Value = 64550.24
$ws.Cells.Item(140, 12).Value = 64550.24
$ws.Cells.Item(140, 14).Value = -74910.23999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15477.932
$ws.Cells.Item(32, 9).Value = 16672.055
$ws.Cells.Item(32, 10).Value = 9666.532999999999
$ws.Cells.Item(32, 11).Value = 16672.055
$ws.Cells.Item(32, 12).Value = 9666.532999999999
$ws.Cells.Item(32, 13).Value = -16385.055
$ws.Cells.Item(32, 14).Value = -10240.533

$ws.Cells.Item(63, 8).Value = 2679
$ws.Cells.Item(63, 9).Value = 2493.9
$ws.Cells.Item(63, 10).Value = 3296
$ws.Cells.Item(63, 11).Value = 2493.9
$ws.Cells.Item(63, 12).Value = 3296
$ws.Cells.Item(63, 13).Value = -1807.9
$ws.Cells.Item(63, 14).Value = -4668

$ws.Cells.Item(66, 8).Value = 2679
$ws.Cells.Item(66, 9).Value = 2493.9
$ws.Cells.Item(66, 10).Value = 3296
$ws.Cells.Item(66, 11).Value = 12469.5
$ws.Cells.Item(66, 12).Value = 16480
$ws.Cells.Item(66, 13).Value = -9037.5
$ws.Cells.Item(66, 14).Value = -23344

$ws.Cells.Item(74, 8).Value = 7638259
$ws.Cells.Item(74, 9).Value = 12551446
$ws.Cells.Item(74, 10).Value = 79509.46000000001
$ws.Cells.Item(74, 11).Value = 12551446
$ws.Cells.Item(74, 12).Value = 79509.46000000001
$ws.Cells.Item(74, 13).Value = -12550572
$ws.Cells.Item(74, 14).Value = -81257.46000000001

$ws.Cells.Item(77, 8).Value = 7638259
$ws.Cells.Item(77, 9).Value = 12551446
$ws.Cells.Item(77, 10).Value = 79509.46000000001
$ws.Cells.Item(77, 11).Value = 62757230
$ws.Cells.Item(77, 12).Value = 397547.3
$ws.Cells.Item(77, 13).Value = -62752862
$ws.Cells.Item(77, 14).Value = -406283.3

$ws.Cells.Item(88, 8).Value = 6304.4165
$ws.Cells.Item(88, 9).Value = 3719
$ws.Cells.Item(88, 11).Value = 3719
$ws.Cells.Item(88, 13).Value = -3313

$ws.Cells.Item(91, 8).Value = 6304.4165
$ws.Cells.Item(91, 9).Value = 3719
$ws.Cells.Item(91, 11).Value = 3719
$ws.Cells.Item(91, 13).Value = -2315

$ws.Cells.Item(122, 8).Value = 6175264
$ws.Cells.Item(122, 9).Value = 2304.4
$ws.Cells.Item(122, 11).Value = 6913.200000000001
$ws.Cells.Item(122, 13).Value = -4463.200000000001

$ws.Cells.Item(132, 8).Value = 45848.297
$ws.Cells.Item(132, 9).Value = 29292.686
$ws.Cells.Item(132, 11).Value = 87878.058
$ws.Cells.Item(132, 13).Value = -85348.058

$ws.Cells.Item(134, 8).Value = 46502.758
$ws.Cells.Item(134, 10).Value = 46502.758
$ws.Cells.Item(134, 12).Value = 46502.758
$ws.Cells.Item(134, 14).Value = -56642.758

$ws.Cells.Item(135, 8).Value = 42243.375
$ws.Cells.Item(135, 10).Value = 42243.375
$ws.Cells.Item(135, 12).Value = 42243.375
$ws.Cells.Item(135, 14).Value = -52383.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1009.96155
$ws.Cells.Item(20, 9).Value = 832.8333
$ws.Cells.Item(20, 11).Value = 832.8333
$ws.Cells.Item(20, 13).Value = -585.8333

$ws.Cells.Item(33, 8).Value = 4000.6667
$ws.Cells.Item(33, 9).Value = 2000
$ws.Cells.Item(33, 11).Value = 2000
$ws.Cells.Item(33, 13).Value = -1664

$ws.Cells.Item(86, 8).Value = 12931.281
$ws.Cells.Item(86, 9).Value = 17306.8
$ws.Cells.Item(86, 10).Value = 9070.529
$ws.Cells.Item(86, 11).Value = 17306.8
$ws.Cells.Item(86, 12).Value = 9070.529
$ws.Cells.Item(86, 13).Value = -16183.8
$ws.Cells.Item(86, 14).Value = -11316.529

$ws.Cells.Item(89, 8).Value = 12931.281
$ws.Cells.Item(89, 9).Value = 17306.8
$ws.Cells.Item(89, 10).Value = 9070.529
$ws.Cells.Item(89, 11).Value = 86534
$ws.Cells.Item(89, 12).Value = 45352.645
$ws.Cells.Item(89, 13).Value = -80918
$ws.Cells.Item(89, 14).Value = -56584.645

$ws.Cells.Item(94, 8).Value = 927
$ws.Cells.Item(94, 9).Value = 918.63635
$ws.Cells.Item(94, 11).Value = 918.63635
$ws.Cells.Item(94, 13).Value = -467.63635

$ws.Cells.Item(99, 8).Value = 1202.6316
$ws.Cells.Item(99, 9).Value = 1263.3334
$ws.Cells.Item(99, 10).Value = 975
$ws.Cells.Item(99, 11).Value = 1263.3334
$ws.Cells.Item(99, 12).Value = 975
$ws.Cells.Item(99, 13).Value = 234.6666
$ws.Cells.Item(99, 14).Value = -3971

$ws.Cells.Item(134, 8).Value = 4054.9614
$ws.Cells.Item(134, 9).Value = 3413.7646
$ws.Cells.Item(134, 10).Value = 5266.1113
$ws.Cells.Item(134, 11).Value = 10241.2938
$ws.Cells.Item(134, 12).Value = 15798.3339
$ws.Cells.Item(134, 13).Value = -7706.293799999999
$ws.Cells.Item(134, 14).Value = -20868.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 37039172
$ws.Cells.Item(58, 9).Value = 41668620
$ws.Cells.Item(58, 11).Value = 41668620
$ws.Cells.Item(58, 13).Value = -41668417

$ws.Cells.Item(62, 9).Value = 3000
$ws.Cells.Item(62, 11).Value = 3000
$ws.Cells.Item(62, 13).Value = -2376

$ws.Cells.Item(65, 9).Value = 3000
$ws.Cells.Item(65, 11).Value = 15000
$ws.Cells.Item(65, 13).Value = -11880

$ws.Cells.Item(135, 8).Value = 48850
$ws.Cells.Item(135, 10).Value = 48850
$ws.Cells.Item(135, 12).Value = 48850
$ws.Cells.Item(135, 14).Value = -58990

$ws.Cells.Item(136, 8).Value = 37039172
$ws.Cells.Item(136, 9).Value = 41668620
$ws.Cells.Item(136, 11).Value = 125005860
$ws.Cells.Item(136, 13).Value = -125003310

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(20, 8).Value = 2500
$ws.Cells.Item(20, 9).Value = 1500
$ws.Cells.Item(20, 10).Value = 3500
$ws.Cells.Item(20, 11).Value = 4500
$ws.Cells.Item(20, 12).Value = 10500
$ws.Cells.Item(20, 13).Value = -4273
$ws.Cells.Item(20, 14).Value = -10954

$ws.Cells.Item(129, 8).Value = 4904067.5
$ws.Cells.Item(129, 9).Value = 2013.75
$ws.Cells.Item(129, 10).Value = 9261448
$ws.Cells.Item(129, 11).Value = 6041.25
$ws.Cells.Item(129, 12).Value = 27784344
$ws.Cells.Item(129, 13).Value = -1041.25
$ws.Cells.Item(129, 14).Value = -27794344

$ws.Cells.Item(131, 8).Value = 996.23376
$ws.Cells.Item(131, 10).Value = 1041.8857
$ws.Cells.Item(131, 12).Value = 3125.6571
$ws.Cells.Item(131, 14).Value = -13205.6571

$ws.Cells.Item(137, 8).Value = 21311.967
$ws.Cells.Item(137, 9).Value = 951.5
$ws.Cells.Item(137, 10).Value = 34885.61
$ws.Cells.Item(137, 11).Value = 2854.5
$ws.Cells.Item(137, 12).Value = 104656.83
$ws.Cells.Item(137, 13).Value = 2245.5
$ws.Cells.Item(137, 14).Value = -114856.83

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 29068.12
$ws.Cells.Item(70, 9).Value = 45038.12
$ws.Cells.Item(70, 10).Value = 5582.8237
$ws.Cells.Item(70, 11).Value = 45038.12
$ws.Cells.Item(70, 12).Value = 5582.8237
$ws.Cells.Item(70, 13).Value = -44768.12
$ws.Cells.Item(70, 14).Value = -6122.8237

$ws.Cells.Item(73, 8).Value = 29068.12
$ws.Cells.Item(73, 9).Value = 45038.12
$ws.Cells.Item(73, 10).Value = 5582.8237
$ws.Cells.Item(73, 11).Value = 45038.12
$ws.Cells.Item(73, 12).Value = 5582.8237
$ws.Cells.Item(73, 13).Value = -44102.12
$ws.Cells.Item(73, 14).Value = -7454.8237

$ws.Cells.Item(131, 8).Value = 54000
$ws.Cells.Item(131, 10).Value = 54000
$ws.Cells.Item(131, 12).Value = 54000
$ws.Cells.Item(131, 14).Value = -64080

$ws.Cells.Item(133, 8).Value = 101390.5
$ws.Cells.Item(133, 10).Value = 101390.5
$ws.Cells.Item(133, 12).Value = 101390.5
$ws.Cells.Item(133, 14).Value = -111510.5

$ws.Cells.Item(135, 8).Value = 30550
$ws.Cells.Item(135, 10).Value = 30550
$ws.Cells.Item(135, 12).Value = 30550
$ws.Cells.Item(135, 14).Value = -40690

$ws.Cells.Item(140, 8).Value = 98410
$ws.Cells.Item(140, 10).Value = 98410
$ws.Cells.Item(140, 12).Value = 98410
$ws.Cells.Item(140, 14).Value = -108770

$ws.Cells.Item(141, 8).Value = 11809.667
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 11809.667
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 11809.667
$ws.Cells.Item(141, 13).Value = ""
$ws.Cells.Item(141, 14).Value = -22169.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 883.2308
$ws.Cells.Item(9, 9).Value = 174.28572
$ws.Cells.Item(9, 10).Value = 1710.3334
$ws.Cells.Item(9, 11).Value = 174.28572
$ws.Cells.Item(9, 12).Value = 1710.3334
$ws.Cells.Item(9, 13).Value = 49.71428
$ws.Cells.Item(9, 14).Value = -2158.3334

$ws.Cells.Item(136, 8).Value = 159580.77
$ws.Cells.Item(136, 9).Value = 126718.75
$ws.Cells.Item(136, 10).Value = 212160
$ws.Cells.Item(136, 11).Value = 380156.25
$ws.Cells.Item(136, 12).Value = 636480
$ws.Cells.Item(136, 13).Value = -377606.25
$ws.Cells.Item(136, 14).Value = -641580
